$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 29 (shifts existing rows 29-50 down to 30-51)
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with this week's data point
$ws.Range("A29").Value = 1
$ws.Range("B29").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C29").Value = "Arica y Parinacota"
$ws.Range("D29").Value = 44596
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = 100112031
$ws.Range("G29").Value = "Poroto verde"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 1300
$ws.Range("K29").Value = 1400
$ws.Range("L29").Value = 1500
$ws.Range("M29").Value = 1450
$ws.Range("N29").Value = "$/kilo"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 1450
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
